# Collapse the multi-run titles/captions (one <a:r> per word) into a
# single run per paragraph, matching the "merged text runs" golden
# fixture. PowerPoint's COM TextRange.Text setter tries to preserve
# existing runs by diffing old vs. new text, so a same-text assignment
# is a no-op and a near-identical assignment only patches the differing
# character(s). Writing an unrelated placeholder value first forces a
# full run collapse, then the second assignment lands the real text in
# that single run.
function Set-MergedText($shape, [string]$text) {
    $shape.TextFrame.TextRange.Text = "zzz___placeholder___zzz"
    $shape.TextFrame.TextRange.Text = $text
}

$p = $ppt.ActivePresentation

# Slide 1: "Section Header (with background image)" title
$s1 = $p.Slides.Item(1)
Set-MergedText $s1.Shapes.Item(1) "Section Header (with background image)"

# Slide 2: "Slide 1" title
$s2 = $p.Slides.Item(2)
Set-MergedText $s2.Shapes.Item(1) "Slide 1"

# Slide 3: "Slide 2" title
$s3 = $p.Slides.Item(3)
Set-MergedText $s3.Shapes.Item(1) "Slide 2"

# Slide 4: "Slide 3" title
$s4 = $p.Slides.Item(4)
Set-MergedText $s4.Shapes.Item(1) "Slide 3"

# Slide 5: "Slide 4" title and the "An image" caption textbox
$s5 = $p.Slides.Item(5)
Set-MergedText $s5.Shapes.Item(1) "Slide 4"
Set-MergedText $s5.Shapes.Item(4) "An image"

# Slide 6: blank-layout slide whose notes page holds the caption text
$s6 = $p.Slides.Item(6)
$notes = $s6.NotesPage
Set-MergedText $notes.Shapes.Item(2) "Blank slides can have background images."
